# Update column G ("K") values for rows 2-24 on the active sheet.
# These new values reflect the commit's regeneration of the "K" (strikeout)
# column using the new calculation method (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 3
    4  = 3
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 0
    17 = 0
    18 = 2
    19 = 0
    20 = 2
    21 = 4
    22 = 2
    23 = 3
    24 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
